$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.133.11'
$ws.Range('E2').Value = '  +4.44%  '
$ws.Range('D3').Value = '2.250.90'
$ws.Range('E3').Value = '  +3.47%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.06'
$ws.Range('E5').Value = '  +3.07%  '
$ws.Range('E6').Value = '  +1.41%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '75.60'
$ws.Range('E7').Value = '  +8.00%  '
$ws.Range('E8').Value = '  -0.14%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.609'
$ws.Range('E9').Value = '  +7.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.07'
$ws.Range('E10').Value = '  +4.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0935'
$ws.Range('E11').Value = '  +1.58%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.01'
$ws.Range('E12').Value = '  +4.29%  '
$ws.Range('E13').Value = '  +0.79%  '
$ws.Range('D14').Value = '2.590.43'
$ws.Range('E14').Value = '  +3.55%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.67'
$ws.Range('E15').Value = '  +2.68%  '
$ws.Range('D16').Value = '2.249.27'
$ws.Range('E16').Value = '  +4.13%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.803'
$ws.Range('E17').Value = '  +1.36%  '
$ws.Range('D18').Value = '43.042.35'
$ws.Range('E18').Value = '  +4.69%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000105'
$ws.Range('E19').Value = '  +4.69%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.40'
$ws.Range('E20').Value = '  +1.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.00'
$ws.Range('E21').Value = '  +2.29%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.96'
$ws.Range('E22').Value = '  +6.33%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '230.56'
$ws.Range('E23').Value = '  +1.85%  '
$ws.Range('E24').Value = '  +15.17%  '
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.92'
$ws.Range('E26').Value = '  +1.63%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.49'
$ws.Range('E27').Value = '  +0.47%  '
$ws.Range('B28').Value = 'PancakeSwap'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.25'
$ws.Range('E28').Value = '  +2.04%  '
$ws.Range('B29').Value = 'InjectiveProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '38.75'
$ws.Range('E29').Value = '  +27.13%  '
$ws.Range('E30').Value = '  +2.09%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '173.62'
$ws.Range('E31').Value = '  +3.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.32'
$ws.Range('E32').Value = '  +1.96%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0797'
$ws.Range('E33').Value = '  +4.77%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.32'
$ws.Range('E34').Value = '  +4.10%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.122'
$ws.Range('E35').Value = '  +1.49%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.110'
$ws.Range('E36').Value = '  +7.59%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.36'
$ws.Range('E37').Value = '  +6.50%  '
$ws.Range('E38').Value = '  +19.29%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '13.12'
$ws.Range('E39').Value = '  +11.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.13'
$ws.Range('E40').Value = '  +3.24%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.53'
$ws.Range('E41').Value = '  +2.84%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.205'
$ws.Range('E42').Value = '  +7.77%  '
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '105.74'
$ws.Range('E43').Value = '  +8.53%  '
$ws.Range('B44').Value = 'MultiversX'
$ws.Range('C44').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '59.77'
$ws.Range('E44').Value = '  +1.45%  '
$ws.Range('E45').Value = '  +5.64%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.480'
$ws.Range('E46').Value = '  +29.31%  '
$ws.Range('E47').Value = '  +3.11%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.39'
$ws.Range('E48').Value = '  +9.32%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.11'
$ws.Range('E49').Value = '  +2.93%  '
$ws.Range('E50').Value = '  +2.26%  '
$ws.Range('D51').Value = '2.463.12'
$ws.Range('E51').Value = '  +3.54%  '
